# Insert a new row at position 116 (pushes existing rows 116-140 down to 117-141)
# and populate the new row with the latest weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(116).Insert()

$ws.Cells.Item(116, 1).Value = 6
$ws.Cells.Item(116, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(116, 3).Value = "Metropolitana"
$ws.Cells.Item(116, 4).Value = 44543
$ws.Cells.Item(116, 5).Value = 13
$ws.Cells.Item(116, 6).Value = 100112001
$ws.Cells.Item(116, 7).Value = "Berenjena"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 700
$ws.Cells.Item(116, 11).Value = 7000
$ws.Cells.Item(116, 12).Value = 8000
$ws.Cells.Item(116, 13).Value = 7643
$ws.Cells.Item(116, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(116, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(116, 16).Value = 153
$ws.Cells.Item(116, 17).Value = 50
$ws.Cells.Item(116, 18).Value = "Hortaliza"
